$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.679174
$ws.Range("N2").Value = 8.037521999999999
$ws.Range("O2").Value = 0.02942326717729479
$ws.Range("P2").Value = 0.02942326717729479
$ws.Range("Q2").Value = 0.1822284849
$ws.Range("R2").Value = 1.6400563641
$ws.Range("S2").Value = 0.02942326717729479
$ws.Range("T2").Value = 0.02942326717729479

# Row 3 updates
$ws.Range("O3").Value = 0.2465847468531156
$ws.Range("P3").Value = 0.2465847468531155
$ws.Range("S3").Value = 0.2465847468531156
$ws.Range("T3").Value = 0.2465847468531155

# Row 4 updates
$ws.Range("M4").Value = 65.67046766666668
$ws.Range("N4").Value = 197.011403
$ws.Range("O4").Value = 0.7212072511207682
$ws.Range("P4").Value = 0.7212072511207681
$ws.Range("Q4").Value = 4.466686309127779
$ws.Range("R4").Value = 40.20017678215001
$ws.Range("S4").Value = 0.7212072511207682
$ws.Range("T4").Value = 0.7212072511207681

# Row 5 updates
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2535676666666667
$ws.Range("N5").Value = 0.760703
$ws.Range("O5").Value = 0.002784734848821526
$ws.Range("P5").Value = 0.002784734848821525
$ws.Range("Q5").Value = 0.01724682746111111
$ws.Range("R5").Value = 0.15522144715
$ws.Range("S5").Value = 0.002784734848821526
$ws.Range("T5").Value = 0.002784734848821525
